$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "attendance" (row 34) and "tasks" (row 35) label-only rows are being
# removed; the "userRole" row that followed them shifts up to become the
# new row 34, and a brand new "graduationDate" row is appended as row 35.
$ws.Rows(34).Delete()
$ws.Rows(34).Delete()

# New graduationDate row.
$ws.Range("A35").Value = "graduationDate"

# B35 should end up with the exact same cell style as B3/C3 ("02/03/2017"
# / "01/22/2017") - copy the date-quote-prefixed format from C3 first,
# then set the (text) value.
$ws.Range("C3").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("B35").Value = "'01/02/2017"

$ws.Range("C35").Value = "'01/06/2017"

$ws.Range("A1").Select()
